$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update defaults to be consistent with the suggested design.
# Column E = HostStructure, Column B = OpticalFiber, Column C = Coating, Column F = HSMVals
$ws.Range("E2").Value = "199"   # HostStructure Young's Modulus (Gpa)
$ws.Range("C4").Value = "2.02"  # Coating Shear modulus (MPa)
$ws.Range("F4").Value = "8.6"   # HSMVals Thermo-optic coeff (ue/K)
$ws.Range("E5").Value = "0.1"   # HostStructure radius/gap/thickness (um/um/1/m)
$ws.Range("B6").Value = "0.55"  # OpticalFiber Coefficient of thermal expansion (ue/K)
$ws.Range("E6").Value = "16"    # HostStructure Coefficient of thermal expansion (ue/K)
